# House_Layout_PPT.pptx — "Fixed bug and corrected house layout"
#
# Slide 1 had a leftover/duplicate cluster of six shapes (three colored
# "bubble" ellipses plus their three "DACH / UK / US" percentage labels)
# sitting right after the "Yard" room box. They are redundant with an
# equivalent cluster elsewhere on the slide, so remove them.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shapeNamesToRemove = @(
    "Ellipse 6",
    "Ellipse 7",
    "Ellipse 8",
    "Textfeld 9",
    "Textfeld 10",
    "Textfeld 11"
)

foreach ($shapeName in $shapeNamesToRemove) {
    $s.Shapes.Item($shapeName).Delete()
}
